# Trade #48 opened / Trade #19 (leadlag) closed - 2026-02-16 21:30:18
# Updates Summary, leadlag, All Trades and Comparison sheets to reflect the
# newly-closed leadlag trade (#19) and the newly-opened leadlag trade (#48).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet: OVERALL + leadlag strategy roll-up numbers
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")

$summary.Range("C2").Value = 19
$summary.Range("D2").Value = "'63.2%"
$summary.Range("E2").Value = "'+2.5392%"
$summary.Range("F2").Value = "'+0.1336%"

$summary.Range("C3").Value = 36
$summary.Range("D3").Value = "'30.6%"
$summary.Range("E3").Value = "'+2.4939%"
$summary.Range("F3").Value = "'+0.0693%"

# ---------------------------------------------------------------------
# leadlag sheet: close trade #19 (row 18) and append newly-opened #48
# ---------------------------------------------------------------------
$leadlag = $wb.Worksheets.Item("leadlag")

# Trade #19 (row 18) transitions from OPEN -> CLOSED
$leadlag.Range("G18").Value = 68950.391145
$leadlag.Range("H18").Value = "CLOSED"
$leadlag.Range("I18").Value = 0.3886
$leadlag.Range("J18").Value = 3.89
$leadlag.Range("M18").Value = "time_exit_5min"
$leadlag.Range("N18").Value = 5

# New row 38: trade #48, freshly opened
$leadlag.Range("A38").Value = 48
$leadlag.Range("B38").Value = "'2026-02-16"
$leadlag.Range("C38").Value = "21:30:18"
$leadlag.Range("D38").Value = "leadlag"
$leadlag.Range("E38").Value = "UP"
$leadlag.Range("F38").Value = 68717.73
$leadlag.Range("G38").Value = "'"
$leadlag.Range("H38").Value = "OPEN"
$leadlag.Range("I38").Value = 0
$leadlag.Range("J38").Value = 0
$leadlag.Range("K38").Value = 0.75
$leadlag.Range("L38").Value = "Binance leading with 0.115% move"
$leadlag.Range("M38").Value = "'"
$leadlag.Range("N38").Value = 0

# ---------------------------------------------------------------------
# All Trades sheet: append the now-closed trade #19 as its own row
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

$allTrades.Range("A20").Value = 19
$allTrades.Range("B20").Value = "'2026-02-16"
$allTrades.Range("C20").Value = "21:25:12"
$allTrades.Range("D20").Value = "leadlag"
$allTrades.Range("E20").Value = "DOWN"
$allTrades.Range("F20").Value = 69219.38
$allTrades.Range("G20").Value = 68950.391145
$allTrades.Range("H20").Value = "CLOSED"
$allTrades.Range("I20").Value = 0.3886
$allTrades.Range("J20").Value = 3.89
$allTrades.Range("K20").Value = 0.75
$allTrades.Range("L20").Value = "Coinbase leading with -0.081% move"
$allTrades.Range("M20").Value = "time_exit_5min"
$allTrades.Range("N20").Value = 5

# ---------------------------------------------------------------------
# Comparison sheet: leadlag strategy stats
# ---------------------------------------------------------------------
$comparison = $wb.Worksheets.Item("Comparison")

$comparison.Range("B2").Value = 36
$comparison.Range("C2").Value = "'30.6%"
$comparison.Range("D2").Value = "'2.08"
$comparison.Range("E2").Value = "'+0.4370%"
$comparison.Range("G2").Value = "'1.13"
